$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift text overhaul -------------------------------------------------
# Jeremy Shade's Saturday/Sunday shift strings ("11am-6.30pm" / "11am - 6.30pm")
# are replaced with a single, consistent split-shift description that is now
# shared by both cells.
$ws.Range("B8").Value = "11am - 3pm/4pm-7.30pm"
$ws.Range("C8").Value = "11am - 3pm/4pm-7.30pm"

# --- Column widths --------------------------------------------------------
# Columns B (Saturday) and C (Sunday) are widened and unified to the same
# width to accommodate the longer shift text.
$ws.Columns.Item(2).ColumnWidth = 26.666666666666668
$ws.Columns.Item(3).ColumnWidth = 26.666666666666668

# --- Selection -------------------------------------------------------------
# The author's cursor ends up on F11 instead of the old I8:I9 selection.
$ws.Range("F11").Select()
